$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 256428.69
$ws.Range("J17").Value = 276046.22
$ws.Range("L17").Value = 828138.6599999999
$ws.Range("N17").Value = -828474.6599999999
$ws.Range("H18").Value = 12258.2
$ws.Range("I18").Value = 14312.5
$ws.Range("J18").Value = 4041
$ws.Range("K18").Value = 14312.5
$ws.Range("L18").Value = 4041
$ws.Range("M18").Value = -14028.5
$ws.Range("N18").Value = -4609
$ws.Range("H29").Value = 4133.3335
$ws.Range("I29").Value = 467.66666
$ws.Range("J29").Value = 5966.1665
$ws.Range("K29").Value = 1402.99998
$ws.Range("L29").Value = 17898.4995
$ws.Range("M29").Value = -1121.99998
$ws.Range("N29").Value = -18460.4995
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("H96").Value = 994.4167
$ws.Range("I96").Value = 954.5
$ws.Range("K96").Value = 2863.5
$ws.Range("M96").Value = -1490.5
$ws.Range("H112").Value = 1861.6111
$ws.Range("J112").Value = 1976
$ws.Range("L112").Value = 5928
$ws.Range("N112").Value = -8144
$ws.Range("H137").Value = 1275.5
$ws.Range("I137").Value = 1202.3334
$ws.Range("K137").Value = 3607.0002
$ws.Range("M137").Value = -1057.0002
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1329.4667
$ws.Range("I2").Value = 1121.8518
$ws.Range("K2").Value = 1121.8518
$ws.Range("M2").Value = -1008.8518
$ws.Range("H32").Value = 5036.3486
$ws.Range("I32").Value = 4553
$ws.Range("K32").Value = 4553
$ws.Range("M32").Value = -4266
$ws.Range("H61").Value = 1771.15
$ws.Range("I61").Value = 1288.4286
$ws.Range("J61").Value = 2897.5
$ws.Range("K61").Value = 1288.4286
$ws.Range("L61").Value = 2897.5
$ws.Range("M61").Value = -1076.4286
$ws.Range("N61").Value = -3321.5
$ws.Range("H63").Value = 2552.5
$ws.Range("I63").Value = 2552.5
$ws.Range("K63").Value = 2552.5
$ws.Range("M63").Value = -1866.5
$ws.Range("H66").Value = 2552.5
$ws.Range("I66").Value = 2552.5
$ws.Range("K66").Value = 12762.5
$ws.Range("M66").Value = -9330.5
$ws.Range("H74").Value = 2014
$ws.Range("I74").Value = 2014
$ws.Range("K74").Value = 2014
$ws.Range("M74").Value = -1140
$ws.Range("H77").Value = 2014
$ws.Range("I77").Value = 2014
$ws.Range("K77").Value = 10070
$ws.Range("M77").Value = -5702
$ws.Range("H93").Value = 29298.666
$ws.Range("J93").Value = 29298.666
$ws.Range("L93").Value = 29298.666
$ws.Range("N93").Value = -34290.666
$ws.Range("H116").Value = 1329.4667
$ws.Range("I116").Value = 1121.8518
$ws.Range("K116").Value = 1121.8518
$ws.Range("M116").Value = 1172.1482
$ws.Range("H136").Value = 1771.15
$ws.Range("I136").Value = 1288.4286
$ws.Range("J136").Value = 2897.5
$ws.Range("K136").Value = 3865.2858
$ws.Range("L136").Value = 8692.5
$ws.Range("M136").Value = -1315.2858
$ws.Range("N136").Value = -13792.5
$ws.Range("H138").Value = 113475.664
$ws.Range("J138").Value = 113475.664
$ws.Range("L138").Value = 113475.664
$ws.Range("N138").Value = -123755.664

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1329.4667
$ws.Range("I3").Value = 1121.8518
$ws.Range("K3").Value = 1121.8518
$ws.Range("M3").Value = -1007.8518
$ws.Range("H22").Value = 3972.0908
$ws.Range("I22").Value = 3972.0908
$ws.Range("K22").Value = 3972.0908
$ws.Range("M22").Value = -3799.0908
$ws.Range("H35").Value = 35000
$ws.Range("J35").Value = 35000
$ws.Range("L35").Value = 35000
$ws.Range("N35").Value = -35620
$ws.Range("H55").Value = 29317.4
$ws.Range("J55").Value = 29317.4
$ws.Range("L55").Value = 29317.4
$ws.Range("N55").Value = -29863.4
$ws.Range("H93").Value = 30448
$ws.Range("J93").Value = 30448
$ws.Range("L93").Value = 30448
$ws.Range("N93").Value = -34192
$ws.Range("H94").Value = 1989.7273
$ws.Range("I94").Value = 1989.7273
$ws.Range("K94").Value = 1989.7273
$ws.Range("M94").Value = -1538.7273

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 900
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("H134").Value = 1968.7858
$ws.Range("I134").Value = 1351
$ws.Range("K134").Value = 4053
$ws.Range("M134").Value = -1518
$ws.Range("H141").Value = 301665.78
$ws.Range("J141").Value = 301665.78
$ws.Range("L141").Value = 301665.78
$ws.Range("N141").Value = -312025.78
$ws.Range("N37").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 99.59999999999999
$ws.Range("I44").Value = 100.75
$ws.Range("J44").Value = 95
$ws.Range("K44").Value = 302.25
$ws.Range("L44").Value = 285
$ws.Range("M44").Value = 95.75
$ws.Range("N44").Value = -1081
$ws.Range("H52").Value = 1429.6666
$ws.Range("J52").Value = 1429.6666
$ws.Range("L52").Value = 4288.9998
$ws.Range("N52").Value = -4820.9998
$ws.Range("H80").Value = 3710.2
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("H83").Value = 3710.2
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("M83").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 30447.5
$ws.Range("J94").Value = 30447.5
$ws.Range("L94").Value = 30447.5
$ws.Range("N94").Value = -31799.5
$ws.Range("H97").Value = 32394.096
$ws.Range("I97").Value = 44273.2
$ws.Range("J97").Value = 2696.3333
$ws.Range("K97").Value = 44273.2
$ws.Range("L97").Value = 2696.3333
$ws.Range("M97").Value = -43777.2
$ws.Range("N97").Value = -3688.3333
$ws.Range("H107").Value = 100007176
$ws.Range("I107").Value = 285
$ws.Range("K107").Value = 285
$ws.Range("M107").Value = 1635
$ws.Range("H132").Value = 2146.4814
$ws.Range("I132").Value = 2146.4814
$ws.Range("K132").Value = 6439.4442
$ws.Range("M132").Value = -3909.4442

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 24630.8
$ws.Range("I93").Value = 2656.077
$ws.Range("K93").Value = 2656.077
$ws.Range("M93").Value = -1408.077
$ws.Range("H100").Value = 29606.6
$ws.Range("I100").Value = 8624
$ws.Range("J100").Value = 47966.375
$ws.Range("K100").Value = 8624
$ws.Range("L100").Value = 47966.375
$ws.Range("M100").Value = -8083
$ws.Range("N100").Value = -49048.375
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("H132").Value = 3729.4348
$ws.Range("I132").Value = 3587.1177
$ws.Range("J132").Value = 4132.6665
$ws.Range("K132").Value = 10761.3531
$ws.Range("L132").Value = 12397.9995
$ws.Range("M132").Value = -8231.3531
$ws.Range("N132").Value = -17457.9995
$ws.Range("N116").ClearContents()

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 18352.273
$ws.Range("J70").Value = 17663
$ws.Range("L70").Value = 17663
$ws.Range("N70").Value = -18293
$ws.Range("H73").Value = 18352.273
$ws.Range("J73").Value = 17663
$ws.Range("L73").Value = 17663
$ws.Range("N73").Value = -19847
$ws.Range("H81").Value = 4000.4
$ws.Range("I81").Value = 4166.6665
$ws.Range("J81").Value = 3751
$ws.Range("K81").Value = 8333.333000000001
$ws.Range("L81").Value = 7502
$ws.Range("M81").Value = -7272.333000000001
$ws.Range("N81").Value = -9624
$ws.Range("H84").Value = 4000.4
$ws.Range("I84").Value = 4166.6665
$ws.Range("J84").Value = 3751
$ws.Range("K84").Value = 41666.665
$ws.Range("L84").Value = 37510
$ws.Range("M84").Value = -36362.665
$ws.Range("N84").Value = -48118
$ws.Range("H100").Value = 4533.222
$ws.Range("I100").Value = 4533.222
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 9066.444
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -8525.444
$ws.Range("H124").Value = 45444.2
$ws.Range("J124").Value = 45444.2
$ws.Range("L124").Value = 45444.2
$ws.Range("N124").Value = -55264.2
$ws.Range("H132").Value = 1470.8928
$ws.Range("I132").Value = 1314.6
$ws.Range("J132").Value = 1861.625
$ws.Range("K132").Value = 3943.8
$ws.Range("L132").Value = 5584.875
$ws.Range("M132").Value = -1413.8
$ws.Range("N132").Value = -10644.875
$ws.Range("N100").ClearContents()
